# The workbook contains a weekly price log for "Pepino ensalada" at
# "Feria Lagunitas de Puerto Montt". A new weekly observation is inserted
# at row 223 (pushing the existing rows 223:261 down to 224:262), and the
# new row is populated with the latest reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 223 - this shifts rows 223:261 down to 224:262
# and carries the existing column formatting (e.g. the date style on D)
# down with them.
$ws.Rows.Item(223).Insert()

# Populate the newly inserted row 223 with the new observation.
$ws.Range("A223").Value2 = 4
$ws.Range("B223").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C223").Value2 = "Los Lagos"
$ws.Range("D223").Value2 = 44694
$ws.Range("E223").Value2 = 10
$ws.Range("F223").Value2 = 100112043
$ws.Range("G223").Value2 = "Pepino ensalada"
$ws.Range("H223").Value2 = "Sin especificar"
$ws.Range("I223").Value2 = "Primera"
$ws.Range("J223").Value2 = 400
$ws.Range("K223").Value2 = 19000
$ws.Range("L223").Value2 = 20000
$ws.Range("M223").Value2 = 19500
$ws.Range("N223").Value2 = "`$/caja 60 unidades"
$ws.Range("O223").Value2 = "Región de Arica y Parinacota"
$ws.Range("P223").Value2 = 325
$ws.Range("Q223").Value2 = 60
$ws.Range("R223").Value2 = "Hortaliza"
